# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Coco" at row 107, pushing the existing
# rows 107-119 down to 108-120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(107).Insert()

$ws.Cells.Item(107, 1).Value = 10
$ws.Cells.Item(107, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(107, 3).Value = "La Araucanía"
$ws.Cells.Item(107, 4).Value = 45154
$ws.Cells.Item(107, 5).Value = 9
$ws.Cells.Item(107, 6).Value = "Fruta"
$ws.Cells.Item(107, 7).Value = 100108
$ws.Cells.Item(107, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(107, 9).Value = 100108007
$ws.Cells.Item(107, 10).Value = "Coco"
$ws.Cells.Item(107, 11).Value = "Sin especificar"
$ws.Cells.Item(107, 12).Value = "Primera"
$ws.Cells.Item(107, 13).Value = 15
$ws.Cells.Item(107, 14).Value = 36000
$ws.Cells.Item(107, 15).Value = 36000
$ws.Cells.Item(107, 16).Value = 36000
$ws.Cells.Item(107, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(107, 18).Value = "Perú"
$ws.Cells.Item(107, 19).Value = 1800
$ws.Cells.Item(107, 20).Value = 20
